$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a header for the new "Electrode Locations" column, matching the style
# already used by the other header cells (A1/B1): bold, bordered, centered.
$ws.Range("C1").Value = "Electrode Locations"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rebuild the data rows (A2:C72): re-sort the records by electrode location
# from A1 through O15 (alphabetically by the row-letter, then numerically by
# the row-number), and populate the new column C with the parsed electrode
# location (e.g. "A11_monopolar_10V_1kHz.txt" -> "A11").
$rows = @(
  @("A1_monopolar_10V_1kHz.txt", 17.321176, "A1"),
  @("A4_monopolar_10V_1kHz.txt", 15.921344, "A4"),
  @("A5_monopolar_10V_1kHz.txt", 14.227521, "A5"),
  @("A7_monopolar_10V_1kHz.txt", 14.829189, "A7"),
  @("A9_monopolar_10V_1kHz.txt", 13.23714, "A9"),
  @("A11_monopolar_10V_1kHz.txt", 9.971793999999999, "A11"),
  @("A13_monopolar_10V_1kHz.txt", 6.278822, "A13"),
  @("B3_monopolar_10V_1kHz.txt", 15.939304, "B3"),
  @("B15_monopolar_10V_1kHz.txt", 9.958323999999999, "B15"),
  @("C1_monopolar_10V_1kHz.txt", 18.382114, "C1"),
  @("C3_monopolar_10V_1kHz.txt", 16.60362, "C3"),
  @("C5_monopolar_10V_1kHz.txt", 14.471481, "C5"),
  @("C7_monopolar_10V_1kHz.txt", 14.289526, "C7"),
  @("C9_monopolar_10V_1kHz.txt", 13.525574, "C9"),
  @("C11_monopolar_10V_1kHz.txt", 10.981631, "C11"),
  @("C14_monopolar_10V_1kHz.txt", 10.896533, "C14"),
  @("C15_monopolar_10V_1kHz.txt", 9.704955999999999, "C15"),
  @("E2_monopolar_10V_1kHz.txt", 18.259599, "E2"),
  @("E3_monopolar_10V_1kHz.txt", 16.649804, "E3"),
  @("E5_monopolar_10V_1kHz.txt", 15.509541, "E5"),
  @("E7_monopolar_10V_1kHz.txt", 15.441335, "E7"),
  @("E9_monopolar_10V_1kHz.txt", 13.961965, "E9"),
  @("E11_monopolar_10V_1kHz.txt", 14.852281, "E11"),
  @("E13_monopolar_10V_1kHz.txt", 11.040643, "E13"),
  @("E15_monopolar_10V_1kHz.txt", 10.218961, "E15"),
  @("F12_monopolar_10V_1kHz.txt", 11.352168, "F12"),
  @("G1_monopolar_10V_1kHz.txt", 18.848012, "G1"),
  @("G3_monopolar_10V_1kHz.txt", 15.951492, "G3"),
  @("G5_monopolar_10V_1kHz.txt", 16.147772, "G5"),
  @("G7_monopolar_10V_1kHz.txt", 15.550379, "G7"),
  @("G9_monopolar_10V_1kHz.txt", 14.749223, "G9"),
  @("G11_monopolar_10V_1kHz.txt", 14.044925, "G11"),
  @("G13_monopolar_10V_1kHz.txt", 11.492215, "G13"),
  @("G15_monopolar_10V_1kHz.txt", 10.48024, "G15"),
  @("H14_monopolar_10V_1kHz.txt", 11.99831, "H14"),
  @("I1_monopolar_10V_1kHz.txt", 19.097959, "I1"),
  @("I3_monopolar_10V_1kHz.txt", 17.180914, "I3"),
  @("I5_monopolar_10V_1kHz.txt", 15.777662, "I5"),
  @("I6_monopolar_10V_1kHz.txt", 13.576889, "I6"),
  @("I7_monopolar_10V_1kHz.txt", 15.100945, "I7"),
  @("I9_monopolar_10V_1kHz.txt", 14.279691, "I9"),
  @("I11_monopolar_10V_1kHz.txt", 14.436202, "I11"),
  @("I13_monopolar_10V_1kHz.txt", 11.569188, "I13"),
  @("I15_monopolar_10V_1kHz.txt", 11.059672, "I15"),
  @("K1_monopolar_10V_1kHz.txt", 19.98635, "K1"),
  @("K3_monopolar_10V_1kHz.txt", 9.902091, "K3"),
  @("K4_monopolar_10V_1kHz.txt", 12.781506, "K4"),
  @("K5_monopolar_10V_1kHz.txt", 14.876228, "K5"),
  @("K7_monopolar_10V_1kHz.txt", 14.56171, "K7"),
  @("K9_monopolar_10V_1kHz.txt", 14.674603, "K9"),
  @("K12_monopolar_10V_1kHz.txt", 12.754351, "K12"),
  @("K13_monopolar_10V_1kHz.txt", 11.671177, "K13"),
  @("K15_monopolar_10V_1kHz.txt", 10.542887, "K15"),
  @("M1_monopolar_10V_1kHz.txt", 21.026976, "M1"),
  @("M3_monopolar_10V_1kHz.txt", 18.924129, "M3"),
  @("M5_monopolar_10V_1kHz.txt", 16.531138, "M5"),
  @("M7_monopolar_10V_1kHz.txt", 15.016275, "M7"),
  @("M9_monopolar_10V_1kHz.txt", 13.884993, "M9"),
  @("M11_monopolar_10V_1kHz.txt", 16.14991, "M11"),
  @("M12_monopolar_10V_1kHz.txt", 11.977143, "M12"),
  @("M13_monopolar_10V_1kHz.txt", 12.355163, "M13"),
  @("M14_monopolar_10V_1kHz.txt", 10.734891, "M14"),
  @("M15_monopolar_10V_1kHz.txt", 11.536475, "M15"),
  @("O1_monopolar_10V_1kHz.txt", 18.191179, "O1"),
  @("O3_monopolar_10V_1kHz.txt", 18.671616, "O3"),
  @("O5_monopolar_10V_1kHz.txt", 17.216194, "O5"),
  @("O7_monopolar_10V_1kHz.txt", 16.256175, "O7"),
  @("O9_monopolar_10V_1kHz.txt", 15.059893, "O9"),
  @("O11_monopolar_10V_1kHz.txt", 14.611528, "O11"),
  @("O13_monopolar_10V_1kHz.txt", 12.705816, "O13"),
  @("O15_monopolar_10V_1kHz.txt", 13.597201, "O15")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

Write-Host "Applied electrode-location sort and added Electrode Locations column."
